$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    # Force text storage so numeric-looking strings (e.g. "251.80", "0.620")
    # keep their exact original formatting instead of being parsed as numbers.
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

$ws.Range('D2').Value = '41.940.46'
$ws.Range('E2').Value = '  +1.57%  '
$ws.Range('D3').Value = '2.210.03'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('E4').Value = '  -0.21%  '
Set-TextCell 'D5' '251.80'
$ws.Range('E5').Value = '  -1.52%  '
Set-TextCell 'D6' '0.622'
$ws.Range('E6').Value = '  -0.46%  '
Set-TextCell 'D7' '67.59'
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  -0.12%  '
Set-TextCell 'D9' '0.620'
$ws.Range('E9').Value = '  +6.99%  '
Set-TextCell 'D10' '38.69'
$ws.Range('E10').Value = '  +1.97%  '
Set-TextCell 'D11' '59.37'
$ws.Range('E11').Value = '  +2.19%  '
$ws.Range('E12').Value = '  -0.52%  '
Set-TextCell 'D13' '7.02'
$ws.Range('E13').Value = '  -0.90%  '
Set-TextCell 'D14' '0.103'
$ws.Range('E14').Value = '  -0.62%  '
$ws.Range('D15').Value = '2.547.49'
$ws.Range('E15').Value = '  +0.96%  '
Set-TextCell 'D16' '0.875'
$ws.Range('E16').Value = '  +0.58%  '
Set-TextCell 'D17' '14.51'
$ws.Range('E17').Value = '  -0.17%  '
$ws.Range('D18').Value = '2.194.31'
$ws.Range('E18').Value = '  +0.12%  '
$ws.Range('D19').Value = '41.846.22'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').Value = '0.0₃0961'
$ws.Range('E20').Value = '  +0.86%  '
$ws.Range('B21').Value = 'Litecoin'
$ws.Range('C21').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextCell 'D21' '72.33'
$ws.Range('E21').Value = '  +0.67%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextCell 'D22' '6.14'
$ws.Range('E22').Value = '  -1.60%  '
Set-TextCell 'D23' '231.13'
$ws.Range('E23').Value = '  -0.60%  '
Set-TextCell 'D24' '2.01'
$ws.Range('E24').Value = '  -3.69%  '
$ws.Range('E25').Value = '  +1.34%  '
$ws.Range('E26').Value = '  -0.08%  '
Set-TextCell 'D27' '11.18'
$ws.Range('E27').Value = '  -6.61%  '
Set-TextCell 'D28' '2.41'
$ws.Range('E28').Value = '  -4.86%  '
Set-TextCell 'D29' '3.69'
$ws.Range('E29').Value = '  -1.13%  '
Set-TextCell 'D30' '2.15'
$ws.Range('E30').Value = '  -1.98%  '
Set-TextCell 'D31' '166.80'
$ws.Range('E31').Value = '  -1.85%  '
Set-TextCell 'D32' '20.40'
$ws.Range('E32').Value = '  -1.03%  '
$ws.Range('B33').Value = 'Kaspa'
$ws.Range('C33').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextCell 'D33' '0.121'
$ws.Range('E33').Value = '  +0.78%  '
$ws.Range('B34').Value = 'InternetComputer(DFINITY)'
$ws.Range('C34').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell 'D34' '5.87'
$ws.Range('E34').Value = '  +7.19%  '
Set-TextCell 'D35' '0.0782'
$ws.Range('E35').Value = '  +7.47%  '
$ws.Range('B37').Value = 'InjectiveProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextCell 'D37' '25.98'
$ws.Range('E37').Value = '  +2.82%  '
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextCell 'D38' '4.58'
$ws.Range('E38').Value = '  -0.67%  '
Set-TextCell 'D39' '4.07'
$ws.Range('E39').Value = '  +2.28%  '
Set-TextCell 'D40' '0.0307'
$ws.Range('E40').Value = '  +3.08%  '
$ws.Range('E41').Value = '  -0.26%  '
$ws.Range('B42').Value = 'THORChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextCell 'D42' '5.64'
$ws.Range('E42').Value = '  -2.09%  '
$ws.Range('B43').Value = 'FTXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextCell 'D43' '5.15'
$ws.Range('E43').Value = '  +5.40%  '
$ws.Range('B44').Value = 'Celestia'
$ws.Range('C44').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextCell 'D44' '11.99'
$ws.Range('E44').Value = '  -1.84%  '
Set-TextCell 'D45' '61.54'
$ws.Range('E45').Value = '  -4.27%  '
Set-TextCell 'D46' '0.195'
$ws.Range('E46').Value = '  -4.50%  '
$ws.Range('E47').Value = '  -0.52%  '
Set-TextCell 'D48' '0.0994'
$ws.Range('E48').Value = '  -2.26%  '
$ws.Range('E49').Value = '  -0.50%  '
Set-TextCell 'D50' '1.15'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('B51').Value = 'HuobiToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextCell 'D51' '2.84'
$ws.Range('E51').Value = '  +5.04%  '
